$d = $word.ActiveDocument

# The document contains 7 occurrences of a split "<id>...</id>" marker where
# the opening tag "<id>", the bare id text (e.g. "p063v_1"), and the closing
# tag "</id>" live in three separate runs. We merge each trio back into a
# single run (keeping the formatting of the first/opening run) so the run's
# text becomes "<id>p063v_N</id>".
#
# Setting Range.Text to a value identical to the range's current text is a
# no-op in this engine, so for each match we first stamp a unique
# placeholder (forcing the merge/replace) and then do a second pass to swap
# the placeholder for the real final text.

for ($i = 1; $i -le 7; $i++) {
    $target = "<id>p063v_$i</id>"
    $placeholder = "@@PLACEHOLDER_$i@@"

    $r = $d.Content
    $found = $r.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Text = $placeholder
    }

    $r2 = $d.Content
    $found2 = $r2.Find.Execute($placeholder, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $r2.Text = $target
    }
}
